$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.472.62"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.852.45"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'240.99"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("D6").Value = "'0.6307"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  +1.68%  "
$ws.Range("D9").Value = "'0.2938"
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("D10").Value = "'24.71"
$ws.Range("E10").Value = "  +0.40%  "
$ws.Range("D11").Value = "'0.07750"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").Value = "1.864.39"
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("D13").Value = "'5.039"
$ws.Range("E13").Value = "  +1.18%  "
$ws.Range("D14").Value = "'0.6815"
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("E15").Value = "  +4.90%  "
$ws.Range("D16").Value = "'83.78"
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("D17").Value = "2.116.47"
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("D18").Value = "'6.182"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").Value = "29.506.49"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").Value = "'229.84"
$ws.Range("E20").Value = "  +0.69%  "
$ws.Range("D21").Value = "'12.46"
$ws.Range("E21").Value = "  +0.41%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").Value = "'7.464"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").Value = "'156.96"
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").Value = "'0.1384"
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("D27").Value = "'8.406"
$ws.Range("E27").Value = "  +0.62%  "
$ws.Range("D28").Value = "'17.71"
$ws.Range("E28").Value = "  +0.46%  "
$ws.Range("D29").Value = "'1.327"
$ws.Range("D30").Value = "'1.470"
$ws.Range("E30").Value = "  +0.96%  "
$ws.Range("D31").Value = "'0.05683"
$ws.Range("E31").Value = "  +1.13%  "
$ws.Range("D32").Value = "'4.139"
$ws.Range("E32").Value = "  +0.47%  "
$ws.Range("D33").Value = "'4.045"
$ws.Range("E33").Value = "  +0.19%  "
$ws.Range("D34").Value = "'1.853"
$ws.Range("E34").Value = "  +0.95%  "
$ws.Range("E35").Value = "  +1.12%  "
$ws.Range("D36").Value = "'0.7081"
$ws.Range("E36").Value = "  -0.84%  "
$ws.Range("D37").Value = "'2.589"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("D38").Value = "'2.783"
$ws.Range("E38").Value = "  +0.58%  "
$ws.Range("E39").Value = "  -0.65%  "
$ws.Range("D40").Value = "1.222.28"
$ws.Range("E40").Value = "  -1.56%  "
$ws.Range("D41").Value = "'6.557"
$ws.Range("E41").Value = "  +5.61%  "
$ws.Range("D42").Value = "'0.9094"
$ws.Range("E42").Value = "  +0.73%  "
$ws.Range("D43").Value = "'1.001"
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "2.024.70"
$ws.Range("E44").Value = "  +0.67%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'101.60"
$ws.Range("E45").Value = "  -0.18%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'66.50"
$ws.Range("E46").Value = "  +0.75%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "'0.00000000120"
$ws.Range("E47").Value = "  +1.48%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.145"
$ws.Range("E48").Value = "  +0.44%  "
$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D49").Value = "'0.4024"
$ws.Range("E49").Value = "  +0.64%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.049"
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "'1.687"
$ws.Range("E51").Value = "  +0.25%  "

$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
